$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the sourceCity/destinationCity values in row 5 (B5 <-> C5)
$ws.Range("B5").Value = "Bengaluru"
$ws.Range("C5").Value = "Delhi"

# Update the active selection to B5:C5 with active cell B5
$ws.Range("B5:C5").Select()
